$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Collapse the "Bugs" bullet list down to a single "None known..." item
#    - Delete the 4 paragraphs from "Other cars can ..." through
#      "Electric vehicles, track-based vehicles and some tugs"
#    - Replace the text of the remaining first bullet with the new text
# -----------------------------------------------------------------------
$pStart = $null
$pEnd = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($pStart -eq $null -and $t -match "Other cars can") {
        $pStart = $p
    }
    if ($t -match "Electric vehicles, track-based vehicles and some tugs") {
        $pEnd = $p
    }
}
if ($pStart -ne $null -and $pEnd -ne $null) {
    $r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $r.Delete()
}

$d.Content.Find.Execute(
    "Car automatically restarts after stalling in versions <= 350",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "None known, please do report any you might find.", 2)

# -----------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker from the run that begins
#    "If your drivers are up-to-date ..." to the run that begins
#    "Currently, this mod only uses the Logitech SDK ...".
#    These markers are not exposed as a scalar COM property, so the
#    paragraphs are rewritten via InsertXML using their exact original
#    markup (minus/plus the lastRenderedPageBreak element).
# -----------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "If your drivers are up-to-date") {
        $full = $p.Range.Duplicate
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00E02EC6" w:rsidRDefault="00E02EC6" w:rsidP="005C478E"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">If your drivers are up-to-date and everything else works, please post your Windows version and the wheel you have. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $full.InsertXML($xml)
        break
    }
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "Currently, this mod only uses") {
        $full = $p.Range.Duplicate
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00E02EC6" w:rsidRPr="00E02EC6" w:rsidRDefault="00E02EC6" w:rsidP="005C478E"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Currently, this mod only uses the Logitech SDK for racing wheels. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ThrustMaster</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> offers no such SDKs, and uses DirectInput, for those wheels you can still just use x360ce or a similar program. DirectInput support is on the to-do list.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $full.InsertXML($xml)
        break
    }
}
